$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value looks numeric need to be forced to Text
# format so Excel keeps them as strings (matching the source data which
# stores all Price/Volume cells as text), instead of auto-converting them
# into numeric cell values.
$textCells = @("D4", "D5", "D6", "D14", "D20", "D21", "D24", "D25", "D26", "D29", "D30", "D36", "D38", "D43", "D46", "D48", "D51")
foreach ($cellAddr in $textCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = '61.586.22'
$ws.Range("E2").Value = '  -2.08%  '
$ws.Range("D3").Value = '2.881.71'
$ws.Range("E3").Value = '  -2.27%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '566.60'
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("D6").Value = '142.11'
$ws.Range("E6").Value = '  -3.47%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D9").Value = '2.879.91'
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("E10").Value = '  -2.44%  '
$ws.Range("E11").Value = '  -2.00%  '
$ws.Range("E12").Value = '  -1.91%  '
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("D14").Value = '31.59'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D16").Value = '3.359.64'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").Value = '61.556.38'
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '2.895.59'
$ws.Range("E18").Value = '  -1.73%  '
$ws.Range("E19").Value = '  -2.76%  '
$ws.Range("D20").Value = '428.94'
$ws.Range("E20").Value = '  -1.91%  '
$ws.Range("D21").Value = '12.97'
$ws.Range("E21").Value = '  -2.98%  '
$ws.Range("E22").Value = '  -2.06%  '
$ws.Range("E23").Value = '  -3.08%  '
$ws.Range("D24").Value = '78.75'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").Value = '11.84'
$ws.Range("E25").Value = '  +0.56%  '
$ws.Range("D26").Value = '10.06'
$ws.Range("E26").Value = '  -10.53%  '
$ws.Range("E28").Value = '  -5.52%  '
$ws.Range("D29").Value = '0.0000106'
$ws.Range("E29").Value = '  +7.07%  '
$ws.Range("D30").Value = '7.03'
$ws.Range("E30").Value = '  -3.42%  '
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("E32").Value = '  -9.53%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  -3.56%  '
$ws.Range("D36").Value = '0.954'
$ws.Range("E36").Value = '  -3.61%  '
$ws.Range("E37").Value = '  -4.45%  '
$ws.Range("D38").Value = '48.76'
$ws.Range("E38").Value = '  -1.66%  '
$ws.Range("E39").Value = '  -7.07%  '
$ws.Range("E40").Value = '  -5.84%  '
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("D43").Value = '39.20'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("E44").Value = '  -4.93%  '
$ws.Range("D45").Value = '2.677.40'
$ws.Range("E45").Value = '  -0.23%  '
$ws.Range("D46").Value = '131.97'
$ws.Range("E46").Value = '  -2.74%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").Value = '343.91'
$ws.Range("E48").Value = '  -3.62%  '
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").Value = '21.35'
$ws.Range("E51").Value = '  -5.57%  '
